$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9645669291338582
$ws.Range("C2").Value = 0.8419243986254296
$ws.Range("D2").Value = 0.8990825688073394

$ws.Range("B3").Value = 0.3185185185185185
$ws.Range("D3").Value = 0.4387755102040816

$ws.Range("B4").Value = 0.8289269051321928
$ws.Range("C4").Value = 0.8289269051321928
$ws.Range("D4").Value = 0.8289269051321928
$ws.Range("E4").Value = 0.8289269051321928

$ws.Range("B5").Value = 0.6415427238261884
$ws.Range("C5").Value = 0.7734212157061574
$ws.Range("D5").Value = 0.6689290395057105

$ws.Range("B6").Value = 0.9032777331034761
$ws.Range("C6").Value = 0.8289269051321928
$ws.Range("D6").Value = 0.8554142475401564
